$d = $word.ActiveDocument

# --- Step 1: append "与建议" as a new, correctly-formatted run right after
#     the existing "第六章 结论" heading run (and its bookmarkEnd id=0). ---

$p1 = $d.Paragraphs.Item(1)
$full = $p1.Range
$full.End = $full.End - 1          ; # exclude the paragraph mark

# Collapse to the end of the heading text and insert a 3-char placeholder;
# this creates a brand-new <w:r> positioned right after bookmarkEnd id=0.
$ins = $full.Duplicate
$ins.Collapse(0)
$ins.InsertAfter("###")

# The whole "第六章 结论" run shares one uniform rPr, so borrow the
# formatting (and, transiently, the text) of its first 3 characters via
# FormattedText; then overwrite the text in place (this keeps the rPr).
$src = $p1.Range
$src.Start = 0
$src.End = 3
$ins.FormattedText = $src.FormattedText
$ins.Text = "与建议"

# --- Step 2: add a fresh "_GoBack" bookmark collapsed right after the new
#     run. Word will slot it in as the next free id and bump every
#     bookmark that follows it by one -- exactly what the diff shows. ---

$p1 = $d.Paragraphs.Item(1)
$endOfHeading = $p1.Range
$endOfHeading.End = $endOfHeading.End - 1
$endOfHeading.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endOfHeading)

# --- Step 3: remove the old "_GoBack" bookmark that used to live alone in
#     the final paragraph, leaving that paragraph empty. ---

$old = $d.Bookmarks.Item("_GoBack")
$oldStart = $old.Range.Start
if ($oldStart -ne $endOfHeading.Start) {
    $old.Delete()
}
